$d = $word.ActiveDocument

# "Versi" + "on"  ->  "Version"
# A Find/Replace that matches exactly the already-contiguous text merges the
# two runs into one without disturbing the surrounding proofErr markers.
$d.Content.Find.Execute("Version", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Version", 2)

# " 2"  ->  " 1."   (the run holding the version number picks up the period)
$d.Range(8, 9).Text = "1."

# Remove the now-redundant trailing "." run that used to follow the
# "_GoBack" bookmark, without touching the bookmark itself.
$d.Range(10, 11).Delete()
